# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Islas Malvinas" / "Groenlandia" rows (208-211 block, rows 209-210) ---
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Update "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 16 de Julio de 2020 a las 08:07"

# --- Row 6 (India) ---
$ws.Range("B6").Value = 970596
$ws.Range("C6").Value = 427
$ws.Range("D6").Value = 613820
$ws.Range("E6").Value = 331841
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = 24935

# --- Row 67 (Uzbekistan) ---
$ws.Range("B67").Value = 14787
$ws.Range("C67").Value = 206
$ws.Range("E67").Value = 6060
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 72

# --- Row 72 (Kirguistan) ---
$ws.Range("B72").Value = 12498
$ws.Range("C72").Value = 521
$ws.Range("D72").Value = 3735
$ws.Range("E72").Value = 8596
$ws.Range("G72").Value = 7
$ws.Range("H72").Value = 167

# --- Row 89 (Haiti) ---
$ws.Range("B89").Value = 6902
$ws.Range("C89").Value = 71
$ws.Range("D89").Value = 3484
$ws.Range("E89").Value = 3273
$ws.Range("G89").Value = 2
$ws.Range("H89").Value = 145

# --- Row 103 (Tailandia) ---
$ws.Range("B103").Value = 3236
$ws.Range("C103").Value = 4
$ws.Range("D103").Value = 3095
$ws.Range("E103").Value = 83

# --- Row 182 (Monaco) ---
$ws.Range("D182").Value = 98
$ws.Range("E182").Value = 7
